# Actualización automática 2025-08-29 13:50:09
$wb = $excel.ActiveWorkbook

# ---- Sheet "VENTAS POR GRUPO" ----
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M3").Value = 489.11
$wsVentasGrupo.Range("R27").Value = 234.9
$wsVentasGrupo.Range("I45").Value = 213.7
$wsVentasGrupo.Range("M55").Value = "24 de 53"
$wsVentasGrupo.Range("R55").Value = "2 de 53"

# ---- Sheet "VENTA MENSUAL" ----
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F3").Value = 756.41
$wsVentaMensual.Range("F27").Value = 6141.83
$wsVentaMensual.Range("F45").Value = 3334.47
$wsVentaMensual.Range("F55").Value = 101715.46

# ---- Sheet "CUMPLIMIENTO MENSUAL" ----
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D8").Value = 2017.33
$wsCumplimiento.Range("E8").Value = -1017.33
$wsCumplimiento.Range("F8").Value = 2.01733

$wsCumplimiento.Range("D13").Value = 274.05
$wsCumplimiento.Range("E13").Value = -254.05
$wsCumplimiento.Range("F13").Value = 13.7025

$wsCumplimiento.Range("D16").Value = 54118.53
$wsCumplimiento.Range("E16").Value = 1941.169999999998
$wsCumplimiento.Range("F16").Value = 0.9653731646797967

$wsCumplimiento.Range("D19").Value = 101715.46
$wsCumplimiento.Range("E19").Value = 15724.23064517915
$wsCumplimiento.Range("F19").Value = 0.866108037591083
